$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New block appended below the existing "Baseline Multiresolution Histograms"
# confusion-matrix table: a fresh confusion matrix (rows 81-88, no row/column
# labels) with per-row Percentage / False Neg / True Pos columns (J:L), plus
# a header row (80) for those three columns and a trailing Accuracy summary
# row (90).
# ---------------------------------------------------------------------------

# Row 80: header labels for J/K/L only (Percentage, False Neg, True Pos)
$ws.Range("J80").Value = "Percentage"
$ws.Range("K80").Value = "False Neg"
$ws.Range("L80").Value = "True Pos"

# Row 81-88: raw confusion-matrix counts (columns B:I)
$data = @(
    @(0, 1, 0, 0, 18, 0, 0, 1),
    @(0, 12, 2, 1, 1, 1, 2, 1),
    @(0, 0, 17, 0, 0, 0, 3, 0),
    @(0, 4, 11, 2, 0, 2, 1, 0),
    @(0, 1, 0, 0, 19, 0, 0, 0),
    @(0, 7, 3, 0, 9, 0, 0, 1),
    @(0, 3, 2, 0, 0, 0, 15, 0),
    @(0, 4, 0, 0, 4, 0, 3, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 81 + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = 2 + $j   # B=2 .. I=9
        $ws.Cells.Item($r, $col).Value = $rowVals[$j]
    }
}

# Column J: percentage of total for each row (relative to $B$22:$I$29)
$ws.Range("J81:J88").Formula = '=SUM(B81:I81)/SUM($B$22:$I$29)'

# Column K: fraction of the row total that is "off-diagonal" for that row's
# designated diagonal column (B81, C82, D83, E84, F85, G86, H87, I88).
$ws.Range("K81").Formula = '=(SUM(B81:I81) - B81) / SUM(B81:I81)'
$ws.Range("K82").Formula = '=(SUM(B82:I82) - C82) / SUM(B82:I82)'
$ws.Range("K83").Formula = '=(SUM(B83:I83) - D83) / SUM(B83:I83)'
$ws.Range("K84").Formula = '=(SUM(B84:I84) - E84) / SUM(B84:I84)'
$ws.Range("K85").Formula = '=(SUM(B85:I85) - F85) / SUM(B85:I85)'
$ws.Range("K86").Formula = '=(SUM(B86:I86) - G86) / SUM(B86:I86)'
$ws.Range("K87").Formula = '=(SUM(B87:I87) - H87) / SUM(B87:I87)'
$ws.Range("K88").Formula = '=(SUM(B88:I88) - I88) / SUM(B88:I88)'

# Column L: complement of K
$ws.Range("L81:L88").Formula = '=1-K81'

# Row 90: overall accuracy label + formula
$ws.Range("K90").Value = "Accuracy"
$ws.Range("K90").Font.Bold = $true
$ws.Range("L90").Formula = '=(B81+C82+D83+E84+F85+G86+H87+I88) / SUM(B81:I88)'

# Move the visible selection to reflect where editing left off
[void]$ws.Range("M80").Select()
